$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (A, C, D get custom widths) ---------------------------
$ws.Columns.Item(1).ColumnWidth = 16.666666666666664
$ws.Columns.Item(3).ColumnWidth = 17.666666666666664
$ws.Columns.Item(4).ColumnWidth = 18.833333333333332

# --- Clear old formatting on the data block so stale styles don't linger -
$ws.Range("B5").ClearFormats()
$ws.Range("C10:D18").ClearFormats()

# --- Summary formulas (Mean / Minimum / Maximum PMD) ----------------------
$ws.Range("B4").Formula = "=AVERAGE(D9:D19)"
$ws.Range("B5").Formula = "=D19"
$ws.Range("D5").Formula = "=B19"
$ws.Range("B6").Formula = "=D17"
$ws.Range("D6").Formula = "=B17"

# --- New header row (row 8) ------------------------------------------------
$ws.Range("A8").ClearContents()
$ws.Range("B8").Value2 = "Wave Lenght (nm)"
$ws.Range("C8").Value2 = "DGD (ps)"
$ws.Range("D8").Value2 = "PMD"

# --- Data table (rows 9-19): index, wavelength, measured DGD, PMD formula -
$waves = @(1550, 1551, 1552, 1553, 1554, 1555, 1556, 1557, 1558, 1559, 1560)
$dgd   = @(2.2316, 2.2778, 2.2865, 2.2169, 2.2389, 2.2442, 2.1282, 2.1219, 2.8889, 2.8124, 1.9269)

for ($i = 0; $i -lt 11; $i++) {
    $r = 9 + $i
    $ws.Cells.Item($r, 1).Value2 = $i + 1
    $ws.Cells.Item($r, 2).Value2 = $waves[$i]
    $ws.Cells.Item($r, 3).Value2 = $dgd[$i]
    $ws.Cells.Item($r, 4).Formula = "=C$r/SQRT(B3)"
}

# --- Number formats ---------------------------------------------------------
$ws.Range("E9:E19").NumberFormat = "0.00"

# --- Selection / view ------------------------------------------------------
$ws.Range("F5").Select()

# --- Page setup --------------------------------------------------------------
$ws.PageSetup.Orientation = 1

Write-Host "applied"
